# Updated cryptos list on Fri Oct 13 15:46:21 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for the coinranking
# snapshot, and swaps ranks 50/51 (Algorand <-> BabyDogeCoin) to reflect
# the new ordering.
#
# Some Price values are plain decimals (e.g. "207.12") that Excel would
# otherwise auto-coerce to a number; the sheet stores these as text, so
# we briefly force a text NumberFormat before writing and restore the
# default "Normal" style afterwards (leaves the cell's effective style
# untouched, only its stored content becomes text like the rest of the
# column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.913.92"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.550.90"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.485"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.772.25"
$ws.Range("D13").Value = "1.547.07"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "26.922.80"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "1.405.14"
$ws.Range("E33").Value = "  +4.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "1.686.28"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0985"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0953"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
